# Auto update Excel log
# Appends newly-logged sensor/alert events to four worksheets:
#   ALERTS     -> 1 new row  (row 10)
#   Proximity  -> 3 new rows (rows 28-30)
#   mmWave     -> 1 new row  (row 10)
#   Camera     -> 2 new rows (rows 16-17)
#
# Helper: writes a row of plain-text values starting at column A of the
# given row. Column A holds date-like text (e.g. "2026-02-01") which Excel
# would otherwise auto-convert into a real date serial; it's temporarily
# forced to Text format, assigned, then the formatting is cleared back off
# so the cell keeps the default "General" style while its stored value
# remains the literal text.
# NOTE: positional parameters are used throughout (named parameters are
# not reliably bound by this host's PowerShell interpreter).
function Write-LogRow {
    param($Sheet, $Row, $Values)

    for ($i = 0; $i -lt $Values.Length; $i++) {
        $cell = $Sheet.Cells.Item($Row, $i + 1)
        if ($i -eq 0) {
            $cell.NumberFormat = "@"
            $cell.Value = $Values[$i]
            $cell.ClearFormats()
        } else {
            $cell.Value = $Values[$i]
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---- ALERTS sheet: new fall-detection alert ----
$wsAlerts = $wb.Worksheets.Item("ALERTS")
Write-LogRow $wsAlerts 10 @("2026-02-01", "14:39:43", "14:00", "Living Room", "CRITICAL", "FALL_DETECTED")

# ---- Proximity sheet: door enter/exit events ----
$wsProximity = $wb.Worksheets.Item("Proximity")
Write-LogRow $wsProximity 28 @("2026-02-01", "14:39:18", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
Write-LogRow $wsProximity 29 @("2026-02-01", "14:39:23", "14:00", "Living Room Main Door", "EXIT", "User EXITED Living Room Main Door")
Write-LogRow $wsProximity 30 @("2026-02-01", "14:39:24", "14:00", "Living Room Main Door", "EXIT", "User EXITED Living Room Main Door")

# ---- mmWave sheet: presence detected ----
$wsMmWave = $wb.Worksheets.Item("mmWave")
Write-LogRow $wsMmWave 10 @("2026-02-01", "14:39:19", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")

# ---- Camera sheet: image captured/received ----
$wsCamera = $wb.Worksheets.Item("Camera")
Write-LogRow $wsCamera 16 @("2026-02-01", "14:39:23", "14:00", "Living Room Main Door", "Image Captured", "Active")
Write-LogRow $wsCamera 17 @("2026-02-01", "14:39:23", "14:00", "Living Room Main Door", "Image Received", "Active")
